$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-normalised LR-pair data (Cort -> Sstr4), including the new
# "Resolving-Mac" sending-cluster block (rows 14-19).
$data = @(
    @("FAPs", "Cort", "Sstr4", "ECs", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 3, 1, 2.935047, 8.805140999999999, 0.5929059285131875, 0.5929059285131874, 1.413284809789, 12.719563288101, 0.4618733921844167, 0.4618733921844166),
    @("FAPs", "Cort", "Sstr4", "FAPs", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 3, 1, 0.448478, 1.345434, 0.09059659521899899, 0.09059659521899897, 0.2159512760526667, 1.943561484474, 0.070574697843027, 0.07057469784302699),
    @("FAPs", "Cort", "Sstr4", "Inflammatory-Mac", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 3, 1, 0.5792206666666667, 1.737662, 0.1170077914200446, 0.1170077914200445, 0.278906528486889, 2.510158756382, 0.09114900515618751, 0.0911490051561875),
    @("FAPs", "Cort", "Sstr4", "MuSCs", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 3, 1, 0.08668166666666666, 0.260045, 0.01751047736546318, 0.01751047736546318, 0.04173898502722222, 0.375650865245, 0.01364065223607398, 0.01364065223607397),
    @("FAPs", "Cort", "Sstr4", "Neutrophils", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 3, 1, 0.6319776666666667, 1.895933, 0.1276651805761876, 0.1276651805761876, 0.3043100967125556, 2.738790870413, 0.09945110544673592, 0.09945110544673591),
    @("FAPs", "Cort", "Sstr4", "Resolving-Mac", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 3, 1, 0.2688693333333333, 0.806608, 0.05431402690611828, 0.05431402690611827, 0.1294660510097778, 1.165194459088, 0.04231059708448599, 0.04231059708448598),
    @("MuSCs", "Cort", "Sstr4", "ECs", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 3, 1, 2.935047, 8.805140999999999, 0.5929059285131875, 0.5929059285131874, 0.161119405065, 1.450074645585, 0.05265518008023936, 0.05265518008023935),
    @("MuSCs", "Cort", "Sstr4", "FAPs", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 3, 1, 0.448478, 1.345434, 0.09059659521899899, 0.09059659521899897, 0.02461919981, 0.22157279829, 0.008045762078776112, 0.00804576207877611),
    @("MuSCs", "Cort", "Sstr4", "Inflammatory-Mac", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 3, 1, 0.5792206666666667, 1.737662, 0.1170077914200446, 0.1170077914200445, 0.03179631849666667, 0.28616686647, 0.01039130498064584, 0.01039130498064584),
    @("MuSCs", "Cort", "Sstr4", "MuSCs", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 3, 1, 0.08668166666666666, 0.260045, 0.01751047736546318, 0.01751047736546318, 0.004758390091666666, 0.042825510825, 0.001555082003112255, 0.001555082003112255),
    @("MuSCs", "Cort", "Sstr4", "Neutrophils", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 3, 1, 0.6319776666666667, 1.895933, 0.1276651805761876, 0.1276651805761876, 0.03469241401166667, 0.312231726105, 0.01133777341385771, 0.0113377734138577),
    @("MuSCs", "Cort", "Sstr4", "Resolving-Mac", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 3, 1, 0.2688693333333333, 0.806608, 0.05431402690611828, 0.05431402690611827, 0.01475958205333333, 0.13283623848, 0.004823555862894384, 0.004823555862894383),
    @("Resolving-Mac", "Cort", "Sstr4", "ECs", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 3, 1, 2.935047, 8.805140999999999, 0.5929059285131875, 0.5929059285131874, 0.239826603766, 2.158439433894, 0.07837735624853141, 0.0783773562485314),
    @("Resolving-Mac", "Cort", "Sstr4", "FAPs", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 3, 1, 0.448478, 1.345434, 0.09059659521899899, 0.09059659521899897, 0.03664573535066666, 0.329811618156, 0.01197613529719588, 0.01197613529719588),
    @("Resolving-Mac", "Cort", "Sstr4", "Inflammatory-Mac", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 3, 1, 0.5792206666666667, 1.737662, 0.1170077914200446, 0.1170077914200445, 0.04732889296755556, 0.425960036708, 0.01546748128321121, 0.01546748128321121),
    @("Resolving-Mac", "Cort", "Sstr4", "MuSCs", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 3, 1, 0.08668166666666666, 0.260045, 0.01751047736546318, 0.01751047736546318, 0.007082874558888888, 0.06374587102999998, 0.00231474312627695, 0.00231474312627695),
    @("Resolving-Mac", "Cort", "Sstr4", "Neutrophils", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 3, 1, 0.6319776666666667, 1.895933, 0.1276651805761876, 0.1276651805761876, 0.05163973778022222, 0.4647576400219999, 0.01687630171559398, 0.01687630171559398),
    @("Resolving-Mac", "Cort", "Sstr4", "Resolving-Mac", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 3, 1, 0.2688693333333333, 0.806608, 0.05431402690611828, 0.05431402690611827, 0.02196967171911111, 0.197727045472, 0.007179873958737906, 0.007179873958737904)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
